$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Natmi LR-pairs output, updated per Dr Hou's advice: the analysis now
# includes a third cluster ("ECs") on the sending side as well (previously
# ECs only appeared as a target cluster), so the Sending/Target cluster
# cross-product grows from 2x3 to 3x3 rows, and all of the derived
# expression/specificity statistics are recomputed accordingly.
#
# Columns: A Sending cluster, B Ligand symbol, C Receptor symbol,
# D Target cluster, E Ligand-expressing cells, F Ligand detection rate,
# G Ligand average expression value, H Ligand total expression value,
# I Ligand derived specificity of average expression value,
# J Ligand derived specificity of total expression value,
# K Receptor-expressing cells, L Receptor detection rate,
# M Receptor average expression value, N Receptor total expression value,
# O Receptor derived specificity of average expression value,
# P Receptor derived specificity of total expression value,
# Q Edge average expression weight, R Edge total expression weight,
# S Edge average expression derived specificity,
# T Edge total expression derived specificity
$rows = @(
    @("ECs","Col18a1","Itga5","ECs", 2, 0.6666666666666666, 10.92859066666667, 32.785772, 0.2185558471001832, 0.2185558471001832, 3, 1, 28.72417333333333, 86.17251999999999, 0.4233259107972328, 0.4233259107972328, 313.9147325983822, 2825.23259338544, 0.09252035303374578, 0.0925203530337458),
    @("ECs","Col18a1","Itga5","FAPs", 2, 0.6666666666666666, 10.92859066666667, 32.785772, 0.2185558471001832, 0.2185558471001832, 3, 1, 30.56986233333333, 91.709587, 0.4505269713084062, 0.4505269713084062, 334.0855121773516, 3006.769609596164, 0.09846530385578864, 0.09846530385578865),
    @("ECs","Col18a1","Itga5","sCs", 2, 0.6666666666666666, 10.92859066666667, 32.785772, 0.2185558471001832, 0.2185558471001832, 3, 1, 8.559531999999999, 25.678596, 0.126147117894361, 0.126147117894361, 93.54362152623466, 841.892593736112, 0.02757019021064874, 0.02757019021064874),
    @("FAPs","Col18a1","Itga5","ECs", 3, 1, 21.275312, 63.825936, 0.4254751576214852, 0.4254751576214852, 3, 1, 28.72417333333333, 86.17251999999999, 0.4233259107972328, 0.4233259107972328, 611.1157496087466, 5500.04174647872, 0.1801146586217114, 0.1801146586217114),
    @("FAPs","Col18a1","Itga5","FAPs", 3, 1, 21.275312, 63.825936, 0.4254751576214852, 0.4254751576214852, 3, 1, 30.56986233333333, 91.709587, 0.4505269713084062, 0.4505269713084062, 650.3833589387147, 5853.450230448432, 0.1916880341301745, 0.1916880341301745),
    @("FAPs","Col18a1","Itga5","sCs", 3, 1, 21.275312, 63.825936, 0.4254751576214852, 0.4254751576214852, 3, 1, 8.559531999999999, 25.678596, 0.126147117894361, 0.126147117894361, 182.106713873984, 1638.960424865856, 0.05367246486959931, 0.05367246486959932),
    @("sCs","Col18a1","Itga5","ECs", 3, 1, 17.79975, 53.39925, 0.3559689952783316, 0.3559689952783316, 3, 1, 28.72417333333333, 86.17251999999999, 0.4233259107972328, 0.4233259107972328, 511.2831042899999, 4601.54793861, 0.1506908991417756, 0.1506908991417756),
    @("sCs","Col18a1","Itga5","FAPs", 3, 1, 17.79975, 53.39925, 0.3559689952783316, 0.3559689952783316, 3, 1, 30.56986233333333, 91.709587, 0.4505269713084062, 0.4505269713084062, 544.1359070677499, 4897.22316360975, 0.1603736333224431, 0.1603736333224431),
    @("sCs","Col18a1","Itga5","sCs", 3, 1, 17.79975, 53.39925, 0.3559689952783316, 0.3559689952783316, 3, 1, 8.559531999999999, 25.678596, 0.126147117894361, 0.126147117894361, 152.357529717, 1371.217767453, 0.04490446281411292, 0.04490446281411293)
)

$r = 2
foreach ($row in $rows) {
    $c = 1
    foreach ($val in $row) {
        $ws.Cells.Item($r, $c).Value = $val
        $c = $c + 1
    }
    $r = $r + 1
}
